$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels in B1/C1 and D1/E1
$ws.Range("B1").Value = "bedrooms_2"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_1"
$ws.Range("E1").Value = "kitchens_2"

# Swap the data values in columns B/C and D/E for each data row (2-7)
for ($row = 2; $row -le 7; $row++) {
    $bVal = $ws.Cells.Item($row, 2).Value()
    $cVal = $ws.Cells.Item($row, 3).Value()
    $ws.Cells.Item($row, 2).Value = $cVal
    $ws.Cells.Item($row, 3).Value = $bVal

    $dVal = $ws.Cells.Item($row, 4).Value()
    $eVal = $ws.Cells.Item($row, 5).Value()
    $ws.Cells.Item($row, 4).Value = $eVal
    $ws.Cells.Item($row, 5).Value = $dVal
}
